$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 6814327
$ws.Cells.Item(2, 3).Value = 'Slovenia Prva Liga'
$ws.Cells.Item(2, 5).Value = 'NS Mura'
$ws.Cells.Item(2, 6).Value = 'NK Domzale'
$ws.Cells.Item(2, 7).Value = 2
$ws.Cells.Item(2, 8).Value = 3
$ws.Cells.Item(2, 9).Value = 2
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 'A'
$ws.Cells.Item(2, 12).Value = 2
$ws.Cells.Item(2, 13).Value = 3.3
$ws.Cells.Item(2, 14).Value = 3.4
$ws.Cells.Item(2, 15).Value = 1.909
$ws.Cells.Item(2, 16).Value = 3.4
$ws.Cells.Item(2, 17).Value = 3.75
$ws.Cells.Item(2, 18).Value = -0.5
$ws.Cells.Item(2, 19).Value = 1.95
$ws.Cells.Item(2, 20).Value = 1.85
$ws.Cells.Item(2, 21).Value = 2.5
$ws.Cells.Item(2, 22).Value = 1.9
$ws.Cells.Item(2, 23).Value = 1.9
$ws.Cells.Item(2, 24).Value = -1
$ws.Cells.Item(2, 25).Value = -1
$ws.Cells.Item(2, 26).Value = 2.75
$ws.Cells.Item(2, 27).Value = -1
$ws.Cells.Item(2, 28).Value = 0.8500000000000001
$ws.Cells.Item(2, 29).Value = 0.8999999999999999
$ws.Cells.Item(2, 30).Value = -1

# Row 3
$ws.Cells.Item(3, 2).Value = 6816473
$ws.Cells.Item(3, 3).Value = 'Slovenia Prva Liga'
$ws.Cells.Item(3, 5).Value = 'NK Bravo'
$ws.Cells.Item(3, 6).Value = 'NK Rogaska'
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 'H'
$ws.Cells.Item(3, 12).Value = 1.8
$ws.Cells.Item(3, 13).Value = 3.5
$ws.Cells.Item(3, 14).Value = 4
$ws.Cells.Item(3, 15).Value = 2.05
$ws.Cells.Item(3, 16).Value = 3
$ws.Cells.Item(3, 17).Value = 3.75
$ws.Cells.Item(3, 18).Value = -0.25
$ws.Cells.Item(3, 19).Value = 1.75
$ws.Cells.Item(3, 20).Value = 2.05
$ws.Cells.Item(3, 21).Value = 2.25
$ws.Cells.Item(3, 22).Value = 1.95
$ws.Cells.Item(3, 23).Value = 1.85
$ws.Cells.Item(3, 24).Value = 1.05
$ws.Cells.Item(3, 25).Value = -1
$ws.Cells.Item(3, 26).Value = -1
$ws.Cells.Item(3, 27).Value = 0.75
$ws.Cells.Item(3, 28).Value = -1
$ws.Cells.Item(3, 29).Value = -0.5
$ws.Cells.Item(3, 30).Value = 0.425

# Row 174
$ws.Cells.Item(174, 2).Value = 7124153
$ws.Cells.Item(174, 3).Value = 'Slovenia Prva Liga'
$ws.Cells.Item(174, 5).Value = 'NK Aluminij'
$ws.Cells.Item(174, 6).Value = 'NK Domzale'
$ws.Cells.Item(174, 7).Value = 1
$ws.Cells.Item(174, 8).Value = 3
$ws.Cells.Item(174, 9).Value = 0
$ws.Cells.Item(174, 10).Value = 3
$ws.Cells.Item(174, 11).Value = 'A'
$ws.Cells.Item(174, 12).Value = 2
$ws.Cells.Item(174, 13).Value = 3.6
$ws.Cells.Item(174, 14).Value = 3
$ws.Cells.Item(174, 15).Value = 1.333
$ws.Cells.Item(174, 16).Value = 4.75
$ws.Cells.Item(174, 17).Value = 7
$ws.Cells.Item(174, 18).Value = -1.5
$ws.Cells.Item(174, 19).Value = 1.95
$ws.Cells.Item(174, 20).Value = 1.85
$ws.Cells.Item(174, 21).Value = 3.25
$ws.Cells.Item(174, 22).Value = 1.95
$ws.Cells.Item(174, 23).Value = 1.85
$ws.Cells.Item(174, 24).Value = -1
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = 6
$ws.Cells.Item(174, 27).Value = -1
$ws.Cells.Item(174, 28).Value = 0.8500000000000001
$ws.Cells.Item(174, 29).Value = 0.95
$ws.Cells.Item(174, 30).Value = -1

# Row 175
$ws.Cells.Item(175, 2).Value = 7124152
$ws.Cells.Item(175, 3).Value = 'Slovenia Prva Liga'
$ws.Cells.Item(175, 5).Value = 'NS Mura'
$ws.Cells.Item(175, 6).Value = 'NK Rogaska'
$ws.Cells.Item(175, 7).Value = 1
$ws.Cells.Item(175, 8).Value = 2
$ws.Cells.Item(175, 9).Value = 0
$ws.Cells.Item(175, 10).Value = 2
$ws.Cells.Item(175, 11).Value = 'A'
$ws.Cells.Item(175, 12).Value = 2.45
$ws.Cells.Item(175, 13).Value = 3.4
$ws.Cells.Item(175, 14).Value = 2.45
$ws.Cells.Item(175, 15).Value = 3.8
$ws.Cells.Item(175, 16).Value = 3.6
$ws.Cells.Item(175, 17).Value = 1.8
$ws.Cells.Item(175, 18).Value = 0.5
$ws.Cells.Item(175, 19).Value = 1.975
$ws.Cells.Item(175, 20).Value = 1.825
$ws.Cells.Item(175, 21).Value = 2.5
$ws.Cells.Item(175, 22).Value = 1.8
$ws.Cells.Item(175, 23).Value = 2
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = -1
$ws.Cells.Item(175, 26).Value = 0.8
$ws.Cells.Item(175, 27).Value = -1
$ws.Cells.Item(175, 28).Value = 0.825
$ws.Cells.Item(175, 29).Value = 0.8
$ws.Cells.Item(175, 30).Value = -1

# Row 176
$ws.Cells.Item(176, 2).Value = 7133777
$ws.Cells.Item(176, 3).Value = 'Slovenia Prva Liga'
$ws.Cells.Item(176, 5).Value = 'NK Radomlje'
$ws.Cells.Item(176, 6).Value = 'NK Celje'
$ws.Cells.Item(176, 7).Value = 1
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 1
$ws.Cells.Item(176, 10).Value = 0
$ws.Cells.Item(176, 11).Value = 'D'
$ws.Cells.Item(176, 12).Value = 3.05
$ws.Cells.Item(176, 13).Value = 3.5
$ws.Cells.Item(176, 14).Value = 2
$ws.Cells.Item(176, 15).Value = 2.9
$ws.Cells.Item(176, 16).Value = 3.6
$ws.Cells.Item(176, 17).Value = 2.1
$ws.Cells.Item(176, 18).Value = 0.25
$ws.Cells.Item(176, 19).Value = 1.9
$ws.Cells.Item(176, 20).Value = 1.9
$ws.Cells.Item(176, 21).Value = 2.75
$ws.Cells.Item(176, 22).Value = 1.8
$ws.Cells.Item(176, 23).Value = 2
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 2.6
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.45
$ws.Cells.Item(176, 28).Value = -0.5
$ws.Cells.Item(176, 29).Value = -1
$ws.Cells.Item(176, 30).Value = 1
